$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.619.49'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '3.447.55'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '581.54'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').Value = '149.29'
$ws.Range('E6').Value = '  +9.33%  '
$ws.Range('D7').Value = '3.448.88'
$ws.Range('E7').Value = '  +2.48%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('E11').Value = '  +3.51%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').Value = '4.035.08'
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').Value = '27.89'
$ws.Range('E14').Value = '  +8.01%  '
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '0.0000175'
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').Value = '3.444.39'
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('D18').Value = '61.719.78'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = '6.29'
$ws.Range('E19').Value = '  +8.48%  '
$ws.Range('D20').Value = '14.38'
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').Value = '389.19'
$ws.Range('E22').Value = '  +4.13%  '
$ws.Range('D23').Value = '0.567'
$ws.Range('E23').Value = '  +2.82%  '
$ws.Range('D24').Value = '3.589.08'
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('D25').Value = '73.10'
$ws.Range('E25').Value = '  +3.01%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').Value = '5.77'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = '0.0000125'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').Value = '7.74'
$ws.Range('E30').Value = '  +3.93%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = '1.54'
$ws.Range('E32').Value = '  -13.55%  '
$ws.Range('D33').Value = '8.25'
$ws.Range('E33').Value = '  +1.76%  '
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').Value = '24.08'
$ws.Range('E36').Value = '  +1.76%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.473.12'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '7.03'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '5.22'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '1.57'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '165.94'
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '0.0787'
$ws.Range('E42').Value = '  +3.38%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '27.10'
$ws.Range('E43').Value = '  +12.23%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.791'
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '4.51'
$ws.Range('E45').Value = '  +2.56%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '42.35'
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '1.71'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.604.69'
$ws.Range('E49').Value = '  +6.19%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').Value = '  -1.90%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '6.96'
$ws.Range('E51').Value = '  +2.51%  '
